$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2_B = New-Object 'object[,]' 1,5
$row2_B[0,0] = 2.345887016753807
$row2_B[0,1] = 0.3135156155690311
$row2_B[0,2] = 0.009519520851597463
$row2_B[0,3] = 0.04641141514002722
$row2_B[0,4] = 4.387006073177275
$ws.Range("B2:F2").Value = $row2_B

$row2_I = New-Object 'object[,]' 1,2
$row2_I[0,0] = 2.666097054353415
$row2_I[0,1] = 0.1207758743605365
$ws.Range("I2:J2").Value = $row2_I

$row2_L = New-Object 'object[,]' 1,2
$row2_L[0,0] = 0.376486380045769
$row2_L[0,1] = 0.5208874751192241
$ws.Range("L2:M2").Value = $row2_L

$row3_B = New-Object 'object[,]' 1,5
$row3_B[0,0] = 2.281277830105239
$row3_B[0,1] = 0.2900822958447122
$row3_B[0,2] = 0.008511352982370823
$row3_B[0,3] = 0.0460698980494314
$row3_B[0,4] = 4.356787725242654
$ws.Range("B3:F3").Value = $row3_B

$row3_I = New-Object 'object[,]' 1,2
$row3_I[0,0] = 2.653196933476536
$row3_I[0,1] = 0.1208226258829952
$ws.Range("I3:J3").Value = $row3_I

$row3_L = New-Object 'object[,]' 1,2
$row3_L[0,0] = 0.3750905753882847
$row3_L[0,1] = 0.5113213601050646
$ws.Range("L3:M3").Value = $row3_L

$row4_B = New-Object 'object[,]' 1,5
$row4_B[0,0] = 2.243105082430418
$row4_B[0,1] = 0.2758938340196266
$row4_B[0,2] = 0.007889797134531307
$row4_B[0,3] = 0.04585534650250533
$row4_B[0,4] = 4.340096582800356
$ws.Range("B4:F4").Value = $row4_B

$row4_I = New-Object 'object[,]' 1,2
$row4_I[0,0] = 2.646304947745122
$row4_I[0,1] = 0.120852234565715
$ws.Range("I4:J4").Value = $row4_I

$row4_L = New-Object 'object[,]' 1,2
$row4_L[0,0] = 0.3744008067487812
$row4_L[0,1] = 0.5057553019091898
$ws.Range("L4:M4").Value = $row4_L

$row5_B = New-Object 'object[,]' 1,5
$row5_B[0,0] = 2.227925915643539
$row5_B[0,1] = 0.2701617096825828
$row5_B[0,2] = 0.007635807815915996
$row5_B[0,3] = 0.04576668139388218
$row5_B[0,4] = 4.333761774287098
$ws.Range("B5:F5").Value = $row5_B

$row5_I = New-Object 'object[,]' 1,2
$row5_I[0,0] = 2.643754160013529
$row5_I[0,1] = 0.1208645252952696
$ws.Range("I5:J5").Value = $row5_I

$row5_L = New-Object 'object[,]' 1,2
$row5_L[0,0] = 0.374161819900003
$row5_L[0,1] = 0.5035644518145546
$ws.Range("L5:M5").Value = $row5_L

$row6_B = New-Object 'object[,]' 1,5
$row6_B[0,0] = 2.225428168320605
$row6_B[0,1] = 0.2692128896897543
$row6_B[0,2] = 0.00759358896642226
$row6_B[0,3] = 0.04575188376765116
$row6_B[0,4] = 4.332738043959196
$ws.Range("B6:F6").Value = $row6_B

$row6_I = New-Object 'object[,]' 1,2
$row6_I[0,0] = 2.643346144090089
$row6_I[0,1] = 0.1208665796854933
$ws.Range("I6:J6").Value = $row6_I

$row6_L = New-Object 'object[,]' 1,2
$row6_L[0,0] = 0.3741246802495013
$row6_L[0,1] = 0.5032053361788869
$ws.Range("L6:M6").Value = $row6_L

$row7_B = New-Object 'object[,]' 1,5
$row7_B[0,0] = 2.242898846328217
$row7_B[0,1] = 0.2758163275451579
$row7_B[0,2] = 0.007886374661932649
$row7_B[0,3] = 0.04585415574393004
$row7_B[0,4] = 4.340009260652707
$ws.Range("B7:F7").Value = $row7_B

$row7_I = New-Object 'object[,]' 1,2
$row7_I[0,0] = 2.646269504556912
$row7_I[0,1] = 0.1208523994145767
$ws.Range("I7:J7").Value = $row7_I

$row7_L = New-Object 'object[,]' 1,2
$row7_L[0,0] = 0.3743974131746768
$row7_L[0,1] = 0.505725442065291
$ws.Range("L7:M7").Value = $row7_L

$row8_B = New-Object 'object[,]' 1,5
$row8_B[0,0] = 2.323298893214599
$row8_B[0,1] = 0.3053940953905396
$row8_B[0,2] = 0.009172393257607325
$row8_B[0,3] = 0.04629466124676007
$row8_B[0,4] = 4.376199182011078
$ws.Range("B8:F8").Value = $row8_B

$row8_I = New-Object 'object[,]' 1,2
$row8_I[0,0] = 2.661435006594388
$row8_I[0,1] = 0.1207918055790747
$ws.Range("I8:J8").Value = $row8_I

$row8_L = New-Object 'object[,]' 1,2
$row8_L[0,0] = 0.3759704138315385
$row8_L[0,1] = 0.5175252480543193
$ws.Range("L8:M8").Value = $row8_L

$row9_B = New-Object 'object[,]' 1,5
$row9_B[0,0] = 2.492860760771862
$row9_B[0,1] = 0.3650031772401405
$row9_B[0,2] = 0.01167704309545314
$row9_B[0,3] = 0.04712050048044603
$row9_B[0,4] = 4.462026099751483
$ws.Range("B9:F9").Value = $row9_B

$row9_I = New-Object 'object[,]' 1,2
$row9_I[0,0] = 2.699382968110882
$row9_I[0,1] = 0.120680238682795
$ws.Range("I9:J9").Value = $row9_I

$row9_L = New-Object 'object[,]' 1,2
$row9_L[0,0] = 0.380381118419379
$row9_L[0,1] = 0.5431064292208845
$ws.Range("L9:M9").Value = $row9_L

$row10_B = New-Object 'object[,]' 1,5
$row10_B[0,0] = 2.624731774308827
$row10_B[0,1] = 0.4098143013108029
$row10_B[0,2] = 0.01351084156436144
$row10_B[0,3] = 0.04770489773864117
$row10_B[0,4] = 4.534258423257029
$ws.Range("B10:F10").Value = $row10_B

$row10_I = New-Object 'object[,]' 1,2
$row10_I[0,0] = 2.732336418376107
$row10_I[0,1] = 0.1206028164658175
$ws.Range("I10:J10").Value = $row10_I

$row10_L = New-Object 'object[,]' 1,2
$row10_L[0,0] = 0.3844296206008551
$row10_L[0,1] = 0.5633946753979444
$ws.Range("L10:M10").Value = $row10_L

$row11_B = New-Object 'object[,]' 1,5
$row11_B[0,0] = 2.686317745346685
$row11_B[0,1] = 0.4304290492365794
$row11_B[0,2] = 0.01434460217723199
$row11_B[0,3] = 0.04796608113073919
$row11_B[0,4] = 4.569137325633534
$ws.Range("B11:F11").Value = $row11_B

$row11_I = New-Object 'object[,]' 1,2
$row11_I[0,0] = 2.748444690373006
$row11_I[0,1] = 0.1205686085019502
$ws.Range("I11:J11").Value = $row11_I

$row11_L = New-Object 'object[,]' 1,2
$row11_L[0,0] = 0.3864468418668565
$row11_L[0,1] = 0.5729501065158047
$ws.Range("L11:M11").Value = $row11_L

$row12_B = New-Object 'object[,]' 1,5
$row12_B[0,0] = 2.709868976687403
$row12_B[0,1] = 0.4382689693157431
$row12_B[0,2] = 0.01466033250911636
$row12_B[0,3] = 0.0480643278329369
$row12_B[0,4] = 4.582637475784765
$ws.Range("B12:F12").Value = $row12_B

$row12_I = New-Object 'object[,]' 1,2
$row12_I[0,0] = 2.754706332899175
$row12_I[0,1] = 0.1205558026803812
$ws.Range("I12:J12").Value = $row12_I

$row12_L = New-Object 'object[,]' 1,2
$row12_L[0,0] = 0.3872359401600249
$row12_L[0,1] = 0.5766154731497579
$ws.Range("L12:M12").Value = $row12_L

$row13_B = New-Object 'object[,]' 1,5
$row13_B[0,0] = 2.70478656505702
$row13_B[0,1] = 0.4365790015305038
$row13_B[0,2] = 0.01459233311432229
$row13_B[0,3] = 0.04804319770106691
$row13_B[0,4] = 4.579716949586441
$ws.Range("B13:F13").Value = $row13_B

$row13_I = New-Object 'object[,]' 1,2
$row13_I[0,0] = 2.753350563159543
$row13_I[0,1] = 0.1205585540206324
$ws.Range("I13:J13").Value = $row13_I

$row13_L = New-Object 'object[,]' 1,2
$row13_L[0,0] = 0.3870648723604688
$row13_L[0,1] = 0.5758239830789336
$ws.Range("L13:M13").Value = $row13_L

$row14_B = New-Object 'object[,]' 1,5
$row14_B[0,0] = 2.68825070877898
$row14_B[0,1] = 0.4310733686542676
$row14_B[0,2] = 0.01437057705108202
$row14_B[0,3] = 0.04797417705650453
$row14_B[0,4] = 4.570242123348493
$ws.Range("B14:F14").Value = $row14_B

$row14_I = New-Object 'object[,]' 1,2
$row14_I[0,0] = 2.74895659073357
$row14_I[0,1] = 0.1205675519791072
$ws.Range("I14:J14").Value = $row14_I

$row14_L = New-Object 'object[,]' 1,2
$row14_L[0,0] = 0.3865112561201443
$row14_L[0,1] = 0.573250717679457
$ws.Range("L14:M14").Value = $row14_L

$row15_B = New-Object 'object[,]' 1,5
$row15_B[0,0] = 2.67815197623446
$row15_B[0,1] = 0.4277053952946517
$row15_B[0,2] = 0.01423474748354892
$row15_B[0,3] = 0.04793181463195939
$row15_B[0,4] = 4.564476633811466
$ws.Range("B15:F15").Value = $row15_B

$row15_I = New-Object 'object[,]' 1,2
$row15_I[0,0] = 2.746286257819804
$row15_I[0,1] = 0.120573082835854
$ws.Range("I15:J15").Value = $row15_I

$row15_L = New-Object 'object[,]' 1,2
$row15_L[0,0] = 0.3861754340229027
$row15_L[0,1] = 0.5716806311213674
$ws.Range("L15:M15").Value = $row15_L

$row16_B = New-Object 'object[,]' 1,5
$row16_B[0,0] = 2.620739149364283
$row16_B[0,1] = 0.4084717483326017
$row16_B[0,2] = 0.01345635014592261
$row16_B[0,3] = 0.04768773605482135
$row16_B[0,4] = 4.532019805758409
$ws.Range("B16:F16").Value = $row16_B

$row16_I = New-Object 'object[,]' 1,2
$row16_I[0,0] = 2.731306285836823
$row16_I[0,1] = 0.1206050726098349
$ws.Range("I16:J16").Value = $row16_I

$row16_L = New-Object 'object[,]' 1,2
$row16_L[0,0] = 0.3843013207995654
$row16_L[0,1] = 0.562776771514983
$ws.Range("L16:M16").Value = $row16_L

$row17_B = New-Object 'object[,]' 1,5
$row17_B[0,0] = 2.585927508899374
$row17_B[0,1] = 0.3967317730724744
$row17_B[0,2] = 0.01297876258028197
$row17_B[0,3] = 0.04753681666201448
$row17_B[0,4] = 4.51262724117268
$ws.Range("B17:F17").Value = $row17_B

$row17_I = New-Object 'object[,]' 1,2
$row17_I[0,0] = 2.722403564886761
$row17_I[0,1] = 0.1206249582103069
$ws.Range("I17:J17").Value = $row17_I

$row17_L = New-Object 'object[,]' 1,2
$row17_L[0,0] = 0.3831965561837194
$row17_L[0,1] = 0.5573981062947908
$ws.Range("L17:M17").Value = $row17_L

$row18_B = New-Object 'object[,]' 1,5
$row18_B[0,0] = 2.566055089040503
$row18_B[0,1] = 0.3900008627710463
$row18_B[0,2] = 0.01270402106481328
$row18_B[0,3] = 0.04744957226350888
$row18_B[0,4] = 4.501663180193304
$ws.Range("B18:F18").Value = $row18_B

$row18_I = New-Object 'object[,]' 1,2
$row18_I[0,0] = 2.717388082092967
$row18_I[0,1] = 0.1206364909647792
$ws.Range("I18:J18").Value = $row18_I

$row18_L = New-Object 'object[,]' 1,2
$row18_L[0,0] = 0.3825776483511873
$row18_L[0,1] = 0.5543351425186245
$ws.Range("L18:M18").Value = $row18_L

$row19_B = New-Object 'object[,]' 1,5
$row19_B[0,0] = 2.55935243776662
$row19_B[0,1] = 0.3877255892724065
$row19_B[0,2] = 0.01261098841030162
$row19_B[0,3] = 0.04741995693376477
$row19_B[0,4] = 4.497983521173126
$ws.Range("B19:F19").Value = $row19_B

$row19_I = New-Object 'object[,]' 1,2
$row19_I[0,0] = 2.715707948592907
$row19_I[0,1] = 0.1206404120124791
$ws.Range("I19:J19").Value = $row19_I

$row19_L = New-Object 'object[,]' 1,2
$row19_L[0,0] = 0.3823709354281988
$row19_L[0,1] = 0.5533033481179217
$ws.Range("L19:M19").Value = $row19_L

$row20_B = New-Object 'object[,]' 1,5
$row20_B[0,0] = 2.589617708472758
$row20_B[0,1] = 0.3979792730550002
$row20_B[0,2] = 0.01302960684514431
$row20_B[0,3] = 0.04755292769052044
$row20_B[0,4] = 4.514671932615812
$ws.Range("B20:F20").Value = $row20_B

$row20_I = New-Object 'object[,]' 1,2
$row20_I[0,0] = 2.723340386696393
$row20_I[0,1] = 0.1206228314955173
$ws.Range("I20:J20").Value = $row20_I

$row20_L = New-Object 'object[,]' 1,2
$row20_L[0,0] = 0.3833124503026539
$row20_L[0,1] = 0.5579674966761061
$ws.Range("L20:M20").Value = $row20_L

$row21_B = New-Object 'object[,]' 1,5
$row21_B[0,0] = 2.693101447044171
$row21_B[0,1] = 0.4326895916016156
$row21_B[0,2] = 0.01443571155672885
$row21_B[0,3] = 0.04799446784389882
$row21_B[0,4] = 4.573017163112809
$ws.Range("B21:F21").Value = $row21_B

$row21_I = New-Object 'object[,]' 1,2
$row21_I[0,0] = 2.750242807345728
$row21_I[0,1] = 0.1205649050244548
$ws.Range("I21:J21").Value = $row21_I

$row21_L = New-Object 'object[,]' 1,2
$row21_L[0,0] = 0.3866731823165992
$row21_L[0,1] = 0.5740052741489379
$ws.Range("L21:M21").Value = $row21_L

$row22_B = New-Object 'object[,]' 1,5
$row22_B[0,0] = 2.762074601580935
$row22_B[0,1] = 0.4555706570091047
$row22_B[0,2] = 0.01535473282848443
$row22_B[0,3] = 0.04827921174054506
$row22_B[0,4] = 4.612853542821284
$ws.Range("B22:F22").Value = $row22_B

$row22_I = New-Object 'object[,]' 1,2
$row22_I[0,0] = 2.768768637299601
$row22_I[0,1] = 0.1205279105167905
$ws.Range("I22:J22").Value = $row22_I

$row22_L = New-Object 'object[,]' 1,2
$row22_L[0,0] = 0.389016610281729
$row22_L[0,1] = 0.5847604620902587
$ws.Range("L22:M22").Value = $row22_L

$row23_B = New-Object 'object[,]' 1,5
$row23_B[0,0] = 2.72513956507845
$row23_B[0,1] = 0.4433405131271115
$row23_B[0,2] = 0.01486420698552848
$row23_B[0,3] = 0.04812758464716271
$row23_B[0,4] = 4.591435569027482
$ws.Range("B23:F23").Value = $row23_B

$row23_I = New-Object 'object[,]' 1,2
$row23_I[0,0] = 2.758794358560621
$row23_I[0,1] = 0.1205475753115248
$ws.Range("I23:J23").Value = $row23_I

$row23_L = New-Object 'object[,]' 1,2
$row23_L[0,0] = 0.3877524351303379
$row23_L[0,1] = 0.5789951740951196
$ws.Range("L23:M23").Value = $row23_L

$row24_B = New-Object 'object[,]' 1,5
$row24_B[0,0] = 2.587948929285574
$row24_B[0,1] = 0.3974152205465771
$row24_B[0,2] = 0.01300662068720726
$row24_B[0,3] = 0.04754564538699046
$row24_B[0,4] = 4.513746951569118
$ws.Range("B24:F24").Value = $row24_B

$row24_I = New-Object 'object[,]' 1,2
$row24_I[0,0] = 2.722916529407513
$row24_I[0,1] = 0.1206237926710814
$ws.Range("I24:J24").Value = $row24_I

$row24_L = New-Object 'object[,]' 1,2
$row24_L[0,0] = 0.3832600039980747
$row24_L[0,1] = 0.5577099840156237
$ws.Range("L24:M24").Value = $row24_L

$row25_B = New-Object 'object[,]' 1,5
$row25_B[0,0] = 2.445711987294203
$row25_B[0,1] = 0.3487015396797517
$row25_B[0,2] = 0.01100088725998205
$row25_B[0,3] = 0.04690107256271414
$row25_B[0,4] = 4.437204897777946
$ws.Range("B25:F25").Value = $row25_B

$row25_I = New-Object 'object[,]' 1,2
$row25_I[0,0] = 2.688231168493871
$row25_I[0,1] = 0.1207096310575482
$ws.Range("I25:J25").Value = $row25_I

$row25_L = New-Object 'object[,]' 1,2
$row25_L[0,0] = 0.3790459967361386
$row25_L[0,1] = 0.5359241212781001
$ws.Range("L25:M25").Value = $row25_L
